$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before the current row 10 (old rows 10-12 shift down
# to 12-14; the existing row 9 stays put and gets filled with new data below).
$ws.Rows.Item(10).Insert()
$ws.Rows.Item(10).Insert()

# Fill row 9 with the new "Baseline 2010 C71" data (mirrors row 8's values).
$ws.Range("A9").Value = "CW3M"
$ws.Range("B9").Value = "Baseline 2010 C71"
$ws.Range("C9").NumberFormat = "0"
$ws.Range("C9").Value = 2010
$ws.Range("D9").Value = 1044.2558590000001
$ws.Range("E9").Value = 1990.4676509999999
$ws.Range("F9").Value = 1.255063
$ws.Range("G9").Value = 327.58108499999997
$ws.Range("H9").Value = 10.610913999999999
$ws.Range("I9").Value = 8.8404570000000007
$ws.Range("J9").Value = 814.39868200000001
$ws.Range("K9").Value = 93.229797000000005
$ws.Range("L9").Value = 1291.7857670000001
$ws.Range("M9").Value = 1165.4420170000001
$ws.Range("N9").Value = 7166.0351559999999
$ws.Range("O9").Value = 29450.638672000001
$ws.Range("P9").Value = -0.473854
$ws.Range("Q9").NumberFormat = "0.000000"
$ws.Range("Q9").Value = -0.00014
$ws.Range("R9").Value = 2010

# The two newly-inserted blank rows (10 and 11) pick up the "C"/"Q" number
# formats used elsewhere in the table, matching the other blank separator row.
$ws.Range("C10").NumberFormat = "0"
$ws.Range("Q10").NumberFormat = "0.000000"
$ws.Range("C11").NumberFormat = "0"
$ws.Range("Q11").NumberFormat = "0.000000"

# Move the selection down to the new blank separator row (old selection was row 9).
$ws.Rows.Item(10).Select() | Out-Null
